$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data cleanup: fix up the pat_id column (A5:A24) which previously was
# a clean sequential 1..20. The refreshed raw export has some duplicate /
# skipped ids (messy_bp source-of-truth regenerated), and the style that
# banded column A in alternating fill colors is removed (reverts to the
# default "Normal" style).
$patIds = @{
    5  = 1
    6  = 2
    7  = 3
    8  = 3
    9  = 4
    10 = 5
    11 = 6
    12 = 7
    13 = 8
    14 = 8
    15 = 11
    16 = 12
    17 = 13
    18 = 14
    19 = 15
    20 = 16
    21 = 17
    22 = 19
    23 = 20
    24 = 21
}

foreach ($row in $patIds.Keys) {
    $cell = $ws.Range("A$row")
    $cell.Value = $patIds[$row]
    $cell.Style = "Normal"
}

# --- Remove the stray "Treatment" / drug-name rows that had been appended
# below the real table (rows 26-28); this also drops the three now-unused
# shared strings ("Treatment", "superpril", "wondersaartan") and shrinks
# the sheet's used range back down to row 24.
$ws.Rows("26:28").Delete()

# --- Selection left where the editor's cursor ended up after the cleanup.
$ws.Range("O15").Select()

Write-Output "edit complete"
